$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 297
$ws.Range("F4").Value = 1230
$ws.Range("F5").Value = 348
$ws.Range("F7").Value = 3815
$ws.Range("F10").Value = 1460
$ws.Range("F14").Value = 144
$ws.Range("F16").Value = 2057
$ws.Range("G17").Value = 55
$ws.Range("F21").Value = 218
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 28
$ws.Range("F4").Value = 36
$ws.Range("F10").Value = 86
$ws.Range("F12").Value = 78
$ws.Range("F13").Value = 224
$ws.Range("F18").Value = 40
$ws.Range("F23").Value = 54
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 810
$ws.Range("F4").Value = 2065
$ws.Range("F5").Value = 302
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 810
$ws.Range("F4").Value = 2065
$ws.Range("F5").Value = 302
$ws.Range("F8").Value = 28
$ws.Range("F9").Value = 36
$ws.Range("F10").Value = 36
$ws.Range("F12").Value = 297
$ws.Range("F13").Value = 1230
$ws.Range("F14").Value = 348
$ws.Range("F19").Value = 3815
$ws.Range("F22").Value = 86
$ws.Range("F24").Value = 78
$ws.Range("F26").Value = 1460
$ws.Range("F28").Value = 224
$ws.Range("F31").Value = 144
$ws.Range("F34").Value = 2057
$ws.Range("G35").Value = 55
$ws.Range("F41").Value = 218
$ws.Range("F44").Value = 40
$ws.Range("F49").Value = 54
